# Weekly refresh of the "Hortaliza, Terminal Hortofrutícola Agro Chillán -
# Repollo" sheet: a new weekly record is inserted above the existing row 47,
# pushing all the later rows down by one (old row 47 -> new row 48, ...,
# old row 187 -> new row 188).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 47..187 down one position, leaving a blank (but pre-formatted,
# inheriting the formatting from the row above) row 47 to fill in.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A47").Value = 7
$ws.Range("B47").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C47").Value = 'Ñuble'
$ws.Range("D47").Value = 44620
$ws.Range("E47").Value = 16
$ws.Range("F47").Value = 100112006
$ws.Range("G47").Value = 'Repollo'
$ws.Range("H47").Value = 'Crespo record'
$ws.Range("I47").Value = 'Primera'
$ws.Range("J47").Value = 300
$ws.Range("K47").Value = 850
$ws.Range("L47").Value = 900
$ws.Range("M47").Value = 875
$ws.Range("N47").Value = '$/unidad'
$ws.Range("O47").Value = 'Provincia de Diguillín'
$ws.Range("P47").Value = 875
$ws.Range("Q47").Value = 1
$ws.Range("R47").Value = 'Hortaliza'
